$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8132768871203098
$ws.Range("C2").Value = 0.2310716155182888
$ws.Range("E2").Value = 0.1284940929411462
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.1346533110046835
$ws.Range("H2").Value = 0.3167417066291094
$ws.Range("I2").Value = 0.2096455108910882
$ws.Range("M2").Value = 0.3331240218226839
$ws.Range("O2").Value = 0.7990675131146645
$ws.Range("B3").Value = 0.7099332356210084
$ws.Range("C3").Value = 0.2080910189378358
$ws.Range("E3").Value = 0.1242556862219288
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.1370821295522546
$ws.Range("H3").Value = 0.3218154056346947
$ws.Range("I3").Value = 0.2161436260102501
$ws.Range("M3").Value = 0.2940286829703851
$ws.Range("O3").Value = 0.8146012699665732
$ws.Range("B4").Value = 0.6462081827346822
$ws.Range("C4").Value = 0.1939196722103134
$ws.Range("E4").Value = 0.1218135239256384
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.1388163646508644
$ws.Range("H4").Value = 0.3251708810746337
$ws.Range("I4").Value = 0.2203860473140824
$ws.Range("M4").Value = 0.2700133019120301
$ws.Range("O4").Value = 0.8251490032327666
$ws.Range("B5").Value = 0.6201733492484323
$ws.Range("C5").Value = 0.1881298082265062
$ws.Range("E5").Value = 0.1208582924457957
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.1395838215030878
$ws.Range("H5").Value = 0.3265985821260209
$ws.Range("I5").Value = 0.2221781865409618
$ws.Range("M5").Value = 0.2602244611590834
$ws.Range("O5").Value = 0.8297003373341312
$ws.Range("B6").Value = 0.6158463347947531
$ws.Range("C6").Value = 0.1871675177706038
$ws.Range("E6").Value = 0.1207020815798252
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.1397149156579438
$ws.Range("H6").Value = 0.326839291375844
$ws.Range("I6").Value = 0.2224795886280666
$ws.Range("M6").Value = 0.2585988971765687
$ws.Range("O6").Value = 0.8304713395413259
$ws.Range("B7").Value = 0.6458573339502323
$ws.Range("C7").Value = 0.1938416479032696
$ws.Range("E7").Value = 0.121800479932773
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.1388264693585377
$ws.Range("H7").Value = 0.3251898914380575
$ws.Range("I7").Value = 0.2204099606076717
$ws.Range("M7").Value = 0.2698812953975818
$ws.Range("O7").Value = 0.8252093607487723
$ws.Range("B8").Value = 0.7777015059635914
$ws.Range("C8").Value = 0.2231608666239708
$ws.Range("E8").Value = 0.1269992371565394
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.1354401404104664
$ws.Range("H8").Value = 0.3184412139559498
$ws.Range("I8").Value = 0.211833531536366
$ws.Range("M8").Value = 0.3196463389937207
$ws.Range("O8").Value = 0.8042133925328869
$ws.Range("B9").Value = 1.034019276943866
$ws.Range("C9").Value = 0.280153201604179
$ws.Range("E9").Value = 0.1384806348225638
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.1307425950755885
$ws.Range("H9").Value = 0.3071166366379856
$ws.Range("I9").Value = 0.1970272064934968
$ws.Range("M9").Value = 0.4171420470466103
$ws.Range("O9").Value = 0.7710945264128384
$ws.Range("B10").Value = 1.220897994007998
$ws.Range("C10").Value = 0.3217000634930685
$ws.Range("E10").Value = 0.147722447925986
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.1284971687539169
$ws.Range("H10").Value = 0.299965465015319
$ws.Range("I10").Value = 0.1873863338510553
$ws.Range("M10").Value = 0.4887109577657895
$ws.Range("O10").Value = 0.7517278980038355
$ws.Range("B11").Value = 1.305586222308023
$ws.Range("C11").Value = 0.3405264204489811
$ws.Range("E11").Value = 0.1521066924078411
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.1277420597659429
$ws.Range("H11").Value = 0.2969670889692395
$ws.Range("I11").Value = 0.1832714078439261
$ws.Range("M11").Value = 0.5212558226177464
$ws.Range("O11").Value = 0.7440076948358012
$ws.Range("B12").Value = 1.337607205788231
$ws.Range("C12").Value = 0.3476445043698391
$ws.Range("E12").Value = 0.1537931558934531
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.1274947868736334
$ws.Range("H12").Value = 0.2958684025641389
$ws.Range("I12").Value = 0.1817523252190028
$ws.Range("M12").Value = 0.5335777678249656
$ws.Range("O12").Value = 0.7412419365277572
$ws.Range("B13").Value = 1.330713108377267
$ws.Range("C13").Value = 0.3461119965879789
$ws.Range("E13").Value = 0.1534287730097716
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.1275463157608101
$ws.Range("H13").Value = 0.2961033889537887
$ws.Range("I13").Value = 0.1820777421716118
$ws.Range("M13").Value = 0.5309241146655808
$ws.Range("O13").Value = 0.7418305616598673
$ws.Range("B14").Value = 1.308221593758731
$ws.Range("C14").Value = 0.3411122536503797
$ws.Range("E14").Value = 0.1522449107895412
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.1277209393487055
$ws.Range("H14").Value = 0.2968759625883663
$ws.Range("I14").Value = 0.1831456461294452
$ws.Range("M14").Value = 0.5222696002391274
$ws.Range("O14").Value = 0.7437769876389098
$ws.Range("B15").Value = 1.294438500759668
$ws.Range("C15").Value = 0.3380483137647445
$ws.Range("E15").Value = 0.1515231892149842
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.1278329483316796
$ws.Range("H15").Value = 0.2973539731177013
$ws.Range("I15").Value = 0.1838048731075266
$ws.Range("M15").Value = 0.5169681783310978
$ws.Range("O15").Value = 0.7449897987183363
$ws.Range("B16").Value = 1.215356727613823
$ws.Range("C16").Value = 0.3204681965325733
$ws.Range("E16").Value = 0.1474395809609916
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.1285519081642192
$ws.Range("H16").Value = 0.3001665506121256
$ws.Range("I16").Value = 0.1876607207280867
$ws.Range("M16").Value = 0.4865838002138645
$ws.Range("O16").Value = 0.7522544466514063
$ws.Range("B17").Value = 1.166758219214785
$ws.Range("C17").Value = 0.3096641933445881
$ws.Range("E17").Value = 0.1449807940596415
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.1290614484932107
$ws.Range("H17").Value = 0.3019572917954321
$ws.Range("I17").Value = 0.1900956592846765
$ws.Range("M17").Value = 0.4679406059395177
$ws.Range("O17").Value = 0.7569909079412724
$ws.Range("B18").Value = 1.13877528243205
$ws.Range("C18").Value = 0.3034431247547502
$ws.Range("E18").Value = 0.143583495715653
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.1293795704180596
$ws.Range("H18").Value = 0.3030112487495273
$ws.Range("I18").Value = 0.1915216434077416
$ws.Range("M18").Value = 0.4572164060132025
$ws.Range("O18").Value = 0.7598177099120278
$ws.Range("B19").Value = 1.129295589022661
$ws.Range("C19").Value = 0.3013356087676584
$ws.Range("E19").Value = 0.143113291806273
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.1294915717106235
$ws.Range("H19").Value = 0.3033722141600421
$ws.Range("I19").Value = 0.1920088253523509
$ws.Range("M19").Value = 0.4535851927870169
$ws.Range("O19").Value = 0.7607923927297549
$ws.Range("B20").Value = 1.17193477041485
$ws.Range("C20").Value = 0.3108150146634898
$ws.Range("E20").Value = 0.1452407813272814
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.1290046120142918
$ws.Range("H20").Value = 0.3017641827806798
$ws.Range("I20").Value = 0.189833818107819
$ws.Range("M20").Value = 0.4699253249504096
$ws.Range("O20").Value = 0.7564760863526629
$ws.Range("B21").Value = 1.314829232276679
$ws.Range("C21").Value = 0.3425811029903798
$ws.Range("E21").Value = 0.1525919247632004
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.1276685957211043
$ws.Range("H21").Value = 0.2966480412201378
$ws.Range("I21").Value = 0.1828309125944454
$ws.Range("M21").Value = 0.5248117002626742
$ws.Range("O21").Value = 0.7432009868029041
$ws.Range("B22").Value = 1.407934748283196
$ws.Range("C22").Value = 0.3632774584614253
$ws.Range("E22").Value = 0.1575494606104684
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.1270210035100874
$ws.Range("H22").Value = 0.2935185061079295
$ws.Range("I22").Value = 0.1784824084543888
$ws.Range("M22").Value = 0.5606707908814883
$ws.Range("O22").Value = 0.7354446900111782
$ws.Range("B23").Value = 1.35826926378644
$ws.Range("C23").Value = 0.3522374932099126
$ws.Range("E23").Value = 0.1548894034230486
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.1273458740224243
$ws.Range("H23").Value = 0.2951691680463071
$ws.Range("I23").Value = 0.1807823261862551
$ws.Range("M23").Value = 0.5415333574522094
$ws.Range("O23").Value = 0.7394998830589827
$ws.Range("B24").Value = 1.169594585425273
$ws.Range("C24").Value = 0.3102947585348375
$ws.Range("E24").Value = 0.1451231903665828
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.1290302293796728
$ws.Range("H24").Value = 0.3018514112574167
$ws.Range("I24").Value = 0.1899521152032762
$ws.Range("M24").Value = 0.4690280520888592
$ws.Range("O24").Value = 0.7567085141033374
$ws.Range("B25").Value = 0.9649253721765945
$ws.Range("C25").Value = 0.2647910519586674
$ws.Range("E25").Value = 0.1352345880743542
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.1318031785338789
$ws.Range("H25").Value = 0.3099752983607189
$ws.Range("I25").Value = 0.2008161963779203
$ws.Range("M25").Value = 0.3907773259919765
$ws.Range("O25").Value = 0.7791860230205856
